# The workbook has a single worksheet ("Taul1") with a small table.
# The edit renames the header of column A from "Name" to "Title",
# and leaves the cursor/selection on cell D6 (as captured in the
# sheetView's <selection> element when the file was last saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the column A header text from "Name" to "Title".
$ws.Range("A1").Value = "Title"

# Update the active selection to match the saved sheet view (D6).
$ws.Range("D6").Select()
